$d = $word.ActiveDocument

$d.TrackRevisions = $true

$rng = $d.Content
$rng.Start = 0
$rng.End = 0

while ($rng.Find.Execute("From the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    # $rng is now collapsed to the found "From the" text.
    # Shrink it to just the leading "F" and retype it together with
    # "Translated " so the run splits into "Translated f" + "rom the".
    $rng.End = $rng.Start + 1
    $rng.Select()
    $word.Selection.TypeText("Translated f")

    # Continue searching after the text we just produced.
    $rng = $d.Content
    $rng.Start = $word.Selection.End
    $rng.End = $d.Content.End
}

$d.TrackRevisions = $false
$d.Revisions.AcceptAll()
